$wb = $excel.ActiveWorkbook

# Add the new "taskManagement" worksheet at the end of the workbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "taskManagement"
$newSheet.Move($wb.Worksheets.Item($wb.Worksheets.Count))

# Header row
$newSheet.Range("A1").Value = "addTask"
$newSheet.Range("B1").Value = "uploadTasks"
$newSheet.Range("C1").Value = "viewTasks"

# Data row
$newSheet.Range("A2").Value = "Add Task"
$newSheet.Range("B2").Value = "Upload Tasks"
$newSheet.Range("C2").Value = "View Tasks"

# Column widths
$newSheet.Columns.Item(1).ColumnWidth = 15.75
$newSheet.Columns.Item(2).ColumnWidth = 16.375
$newSheet.Columns.Item(3).ColumnWidth = 15.625

# Make the new sheet the active/selected tab
$newSheet.Select()
$newSheet.Range("C3").Select()
